$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 44.24836572320267
$ws.Range("R2").Value = 398.235291508824
$ws.Range("S2").Value = 0.006187370219277298
$ws.Range("T2").Value = 0.006187370219277296

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 633.0249328346827
$ws.Range("R3").Value = 5697.224395512144
$ws.Range("S3").Value = 0.08851761084200863
$ws.Range("T3").Value = 0.08851761084200861

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 96.83272190237068
$ws.Range("R4").Value = 871.494497121336
$ws.Range("S4").Value = 0.01354038482456576
$ws.Range("T4").Value = 0.01354038482456575

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 334.6768207122301
$ws.Range("R5").Value = 3012.091386410071
$ws.Range("S5").Value = 0.04679877685225792
$ws.Range("T5").Value = 0.04679877685225792

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 4787.945689971358
$ws.Range("R6").Value = 43091.51120974222
$ws.Range("S6").Value = 0.6695115647652365
$ws.Range("T6").Value = 0.6695115647652364

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 732.4037165558689
$ws.Range("R7").Value = 6591.633449002819
$ws.Range("S7").Value = 0.1024140184668903
$ws.Range("T7").Value = 0.1024140184668903

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 29.85319738807466
$ws.Range("R8").Value = 268.6787764926719
$ws.Range("S8").Value = 0.004174454388319284
$ws.Range("T8").Value = 0.004174454388319283

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 427.0851129215146
$ws.Range("R9").Value = 3843.766016293631
$ws.Range("S9").Value = 0.0597204815499341
$ws.Range("T9").Value = 0.05972048154993409

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 65.33046618397867
$ws.Range("R10").Value = 587.974195655808
$ws.Range("S10").Value = 0.009135338091510304
$ws.Range("T10").Value = 0.009135338091510303
